# Update for release "mines - January 30": refresh the embedded build
# timestamp from "January 30 2026 16.19.47 EST" to "February 02 2026 12.49.33 EST"
# wherever it is stamped into the workbook's text.

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$wb = $excel.ActiveWorkbook

# --- Sheet "About": version banner (A2) and citation text (A6) ---
$wsAbout = $wb.Worksheets.Item("About")

$a2 = $wsAbout.Range("A2").Value()
$wsAbout.Range("A2").Value = $a2.Replace($oldStamp, $newStamp)

$a6 = $wsAbout.Range("A6").Value()
$wsAbout.Range("A6").Value = $a6.Replace($oldStamp, $newStamp)

# --- Sheet "Boundaries and methane sources": build_version column (S2:S9) ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 9; $row++) {
    $cell = $wsData.Cells.Item($row, 19)
    $current = $cell.Value()
    $cell.Value = $current.Replace($oldStamp, $newStamp)
}
